# Colocando header nos gráficos

$wb = $excel.ActiveWorkbook

# --- Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
#     "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
#     Each gets a header label in A1 ("Fonte/Tecnologia") and the A2:A12 labels
#     lose their bold/header style while some get accented spelling fixes.

$sheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Add header label in A1, using the same header style as B1 (copy formats)
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

    # Fix accented labels
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # Remove header style from A2:A12 (now plain / normal style)
    $ws.Range("A2:A12").Style = "Normal"
}

# --- Sheet 5: "Emissoes Totais (MtCO2eq)"
#     Adds header label "Período" in A1, fixes accents in A2/A3, removes row 4 ("Teto").

$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Range("A2:A3").Style = "Normal"

$ws5.Rows.Item(4).Delete()

# --- Sheet 6: "Custo Total (bilhões de R$)"
#     Adds header label "Tipo Expansão" in A1, changes B1 from "Custo" to "2015",
#     fixes accents in A2/A3, updates B2/B3 values.

$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

# B1 text changes from "Custo" to "2015" but must remain a text value (not a
# number), so force text format before assigning, then restore the original
# header style (index 1) by copying formats from an existing header cell.
$wsRef = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
$wsRef.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4122)  # xlPasteFormats

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 168

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99

$ws6.Range("A2:A3").Style = "Normal"
